$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "B" variant row (row 3: B / 10000 / 350) - rows below shift up.
$ws.Rows.Item(3).Delete() | Out-Null

# Relabel the remaining variants (order matters for shared-string table layout):
# row3 (was C) -> H, row4 (was D) -> P, row5 (was E) -> F, row2 (was A) -> Z
$ws.Range("A3").Value2 = "H"
$ws.Range("A4").Value2 = "P"
$ws.Range("A5").Value2 = "F"
$ws.Range("A2").Value2 = "Z"

# Rows 4 & 5 inherited the old "D"/"E" formatting on column A; bring A4/A5 back
# in line with the rest of the label column by copying A3's format onto them.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null

$ws.Range("A4").Select() | Out-Null
